{"js": "// \"Update to lord Barrignton ability\"\n//\n// 1) Remove the stray `_GoBack` bookmark that currently sits in the\n//    \"Time for surgery\" ability line (just leftover cursor-position\n//    metadata from the previous save).\n// 2) In Lord Barrington's 4th ability \"State Titles...\", change\n//    \"they are stunned for 2 turns.\" to \"they fall asleep for 2 turns.\"\n//    and leave the `_GoBack` bookmark positioned right after the newly\n//    typed \"fall asleep\" text (i.e. where the author's cursor ended up),\n//    matching how Word marks the last edit location.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// --- Step 1: drop the old _GoBack bookmark near \"Time for surgery\" ---\nconst oldMark = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\noldMark.load(\"isNullObject\");\nawait context.sync();\nif (!oldMark.isNullObject) {\n  doc.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- Step 2: update the \"State Titles\" ability text ---\nconst target = body.search(\"are stunned\", { matchCase: true, matchWholeWord: false });\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length === 0) {\n  throw new Error('Could not find \"are stunned\" text to update.');\n}\n\nconst hit = target.items[0];\nconst beforePoint = hit.getRange(\"Start\");\nconst afterPoint = hit.getRange(\"End\");\n\n// Drop temporary bookmarks at the exact boundaries of \"are stunned\" so the\n// surrounding text is split into its own runs (mirrors the run layout Word\n// itself produces around an in-place text edit).\nbeforePoint.insertBookmark(\"__tmp_split_start\");\nafterPoint.insertBookmark(\"__tmp_split_end\");\nawait context.sync();\n\n// Re-find \"are stunned\" (still present, now isolated in its own run) and\n// swap it for the new wording.\nconst target2 = body.search(\"are stunned\", { matchCase: true, matchWholeWord: false });\ntarget2.load(\"items\");\nawait context.sync();\ntarget2.items[0].insertText(\"fall asleep\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-plant the real _GoBack bookmark exactly where the temporary \"end\"\n// marker is (right after \"fall asleep\", before \" for 2 turns.\"), then\n// clean up the temporary markers.\nconst endMarker = doc.getBookmarkRange(\"__tmp_split_end\");\nendMarker.insertBookmark(\"_GoBack\");\nawait context.sync();\n\ndoc.deleteBookmark(\"__tmp_split_start\");\ndoc.deleteBookmark(\"__tmp_split_end\");\nawait context.sync();\n", "ps1": "# \"Update to lord Barrignton ability\"\n#\n# 1) Remove the stray `_GoBack` bookmark currently sitting in the\n#    \"Time for surgery\" ability line (leftover cursor-position metadata\n#    from the previous save).\n# 2) In Lord Barrington's 4th ability \"State Titles...\", change\n#    \"they are stunned for 2 turns.\" to \"they fall asleep for 2 turns.\"\n#    and leave the `_GoBack` bookmark positioned right after the newly\n#    typed \"fall asleep\" text, matching where the author's cursor ended\n#    up after the edit.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: drop the old _GoBack bookmark near \"Time for surgery\" ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Step 2: update the \"State Titles\" ability text ---\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"are stunned\")\nif (-not $found) {\n    throw 'Could not find \"are stunned\" text to update.'\n}\n$startPos = $findRange.Start\n$endPos = $findRange.End\n\n# Drop temporary bookmarks at the exact boundaries of \"are stunned\" so the\n# surrounding text splits into its own runs (mirrors the run layout Word\n# itself produces around an in-place text edit) instead of silently\n# re-merging with the neighboring text.\n$startRange = $d.Range($startPos, $startPos)\n$d.Bookmarks.Add(\"__tmp_split_start\", $startRange) | Out-Null\n\n$endRange = $d.Range($endPos, $endPos)\n$d.Bookmarks.Add(\"__tmp_split_end\", $endRange) | Out-Null\n\n# Swap \"are stunned\" for the new wording.\n$replaceRange = $d.Range($startPos, $endPos)\n$replaceRange.Text = \"fall asleep\"\n\n# Re-plant the real _GoBack bookmark exactly where the temporary \"end\"\n# marker is (right after \"fall asleep\", before \" for 2 turns.\"), then\n# clean up the temporary markers.\n$tmpEndBookmark = $d.Bookmarks.Item(\"__tmp_split_end\")\n$d.Bookmarks.Add(\"_GoBack\", $tmpEndBookmark.Range) | Out-Null\n\n$d.Bookmarks.Item(\"__tmp_split_start\").Delete()\n$d.Bookmarks.Item(\"__tmp_split_end\").Delete()\n"}
